$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 53 (data shrinks from 53 to 52 rows)
$ws.Rows.Item(53).Delete()

# Rebuild rows 2:52 with the corrected naive-forecaster data
$data = New-Object 'object[,]' 51,5
$data[0,0] = 39583; $data[0,1] = 2008; $data[0,2] = $null; $data[0,3] = 2009; $data[0,4] = 5.515135932605308
$data[1,0] = 39765; $data[1,1] = 2008; $data[1,2] = $null; $data[1,3] = 2009; $data[1,4] = 2.741420124134053
$data[2,0] = 39948; $data[2,1] = 2009; $data[2,2] = -2.764179109379705; $data[2,3] = 2010; $data[2,4] = -0.2596132895792413
$data[3,0] = 40130; $data[3,1] = 2009; $data[3,2] = -3.872359107260159; $data[3,3] = 2010; $data[3,4] = -1.796159642965267
$data[4,0] = 40310; $data[4,1] = 2010; $data[4,2] = 3.38265053313096; $data[4,3] = 2011; $data[4,4] = 1.421687849828523
$data[5,0] = 40494; $data[5,1] = 2010; $data[5,2] = 4.530477057343663; $data[5,3] = 2011; $data[5,4] = 3.372216600811506
$data[6,0] = 40676; $data[6,1] = 2011; $data[6,2] = 5.813045170083808; $data[6,3] = 2012; $data[6,4] = 2.371160938652705
$data[7,0] = 40862; $data[7,1] = 2011; $data[7,2] = 6.833902841285977; $data[7,3] = 2012; $data[7,4] = 5.794553957309168
$data[8,0] = 41044; $data[8,1] = 2012; $data[8,2] = 4.597183386292891; $data[8,3] = 2013; $data[8,4] = 6.213583554874536
$data[9,0] = 41228; $data[9,1] = 2012; $data[9,2] = 4.166536506645224; $data[9,3] = 2013; $data[9,4] = 4.453225076234824
$data[10,0] = 41409; $data[10,1] = 2013; $data[10,2] = 2.113053977048707; $data[10,3] = 2014; $data[10,4] = 3.80132608443593
$data[11,0] = 41592; $data[11,1] = 2013; $data[11,2] = 2.669880057548091; $data[11,3] = 2014; $data[11,4] = 4.443665867214519
$data[12,0] = 41774; $data[12,1] = 2014; $data[12,2] = 5.262295419893648; $data[12,3] = 2015; $data[12,4] = 3.911207409579687
$data[13,0] = 41957; $data[13,1] = 2014; $data[13,2] = 5.50293301232252; $data[13,3] = 2015; $data[13,4] = 4.93942156312468
$data[14,0] = 42137; $data[14,1] = 2015; $data[14,2] = 4.78031692483154; $data[14,3] = 2016; $data[14,4] = 4.481753591536197
$data[15,0] = 42321; $data[15,1] = 2015; $data[15,2] = 4.829481320500406; $data[15,3] = 2016; $data[15,4] = 5.304093839766777
$data[16,0] = 42503; $data[16,1] = 2016; $data[16,2] = 5.386945580119185; $data[16,3] = 2017; $data[16,4] = 5.291601650505706
$data[17,0] = 42689; $data[17,1] = 2016; $data[17,2] = 5.100281927437122; $data[17,3] = 2017; $data[17,4] = 4.773858406282372
$data[18,0] = 42867; $data[18,1] = 2017; $data[18,2] = 4.884223728030879; $data[18,3] = 2018; $data[18,4] = 4.848884192354119
$data[19,0] = 43053; $data[19,1] = 2017; $data[19,2] = 5.161358932333737; $data[19,3] = 2018; $data[19,4] = 5.448823958443616
$data[20,0] = 43145; $data[20,1] = 2018; $data[20,2] = 6.111775080355519; $data[20,3] = 2019; $data[20,4] = 5.5540247491221
$data[21,0] = 43235; $data[21,1] = 2018; $data[21,2] = 6.072005530313129; $data[21,3] = 2019; $data[21,4] = 5.47980442665561
$data[22,0] = 43326; $data[22,1] = 2018; $data[22,2] = 5.969015988310433; $data[22,3] = 2019; $data[22,4] = 5.321568540215016
$data[23,0] = 43418; $data[23,1] = 2018; $data[23,2] = 5.902681694119694; $data[23,3] = 2019; $data[23,4] = 5.03659417252571
$data[24,0] = 43510; $data[24,1] = 2019; $data[24,2] = 3.94180168325462; $data[24,3] = 2020; $data[24,4] = 5.077291536394846
$data[25,0] = 43600; $data[25,1] = 2019; $data[25,2] = 4.927928448556984; $data[25,3] = 2020; $data[25,4] = 5.802765867180804
$data[26,0] = 43691; $data[26,1] = 2019; $data[26,2] = 3.961704997901161; $data[26,3] = 2020; $data[26,4] = 4.342221082693465
$data[27,0] = 43783; $data[27,1] = 2019; $data[27,2] = 3.884502719230132; $data[27,3] = 2020; $data[27,4] = 4.009670676786059
$data[28,0] = 43875; $data[28,1] = 2020; $data[28,2] = 3.818894565497888; $data[28,3] = 2021; $data[28,4] = 4.451460339319802
$data[29,0] = 43966; $data[29,1] = 2020; $data[29,2] = 2.502519143054571; $data[29,3] = 2021; $data[29,4] = 3.379658261193086
$data[30,0] = 44068; $data[30,1] = 2020; $data[30,2] = -3.840397826549158; $data[30,3] = 2021; $data[30,4] = -2.491095770678031
$data[31,0] = 44159; $data[31,1] = 2020; $data[31,2] = -3.840397826549158; $data[31,3] = 2021; $data[31,4] = -0.236010050592228
$data[32,0] = 44251; $data[32,1] = 2021; $data[32,2] = -1.162988086281536; $data[32,3] = 2022; $data[32,4] = -0.8351862827690737
$data[33,0] = 44341; $data[33,1] = 2021; $data[33,2] = 0.3986977119751156; $data[33,3] = 2022; $data[33,4] = 1.500453746466346
$data[34,0] = 44432; $data[34,1] = 2021; $data[34,2] = 0.4839811651348835; $data[34,3] = 2022; $data[34,4] = 1.25489130894445
$data[35,0] = 44525; $data[35,1] = 2021; $data[35,2] = 0.4839811651348835; $data[35,3] = 2022; $data[35,4] = 2.177145583294293
$data[36,0] = 44617; $data[36,1] = 2022; $data[36,2] = 0.9587999512773893; $data[36,3] = 2023; $data[36,4] = -1.150174078429844
$data[37,0] = 44706; $data[37,1] = 2022; $data[37,2] = 1.97557360987699; $data[37,3] = 2023; $data[37,4] = 0.05919149746531627
$data[38,0] = 44798; $data[38,1] = 2022; $data[38,2] = 2.06342951900429; $data[38,3] = 2023; $data[38,4] = 0.47378583592943
$data[39,0] = 44890; $data[39,1] = 2022; $data[39,2] = 2.06342951900429; $data[39,3] = 2023; $data[39,4] = 0.3896432785800652
$data[40,0] = 44981; $data[40,1] = 2023; $data[40,2] = -2.242084520390608; $data[40,3] = 2024; $data[40,4] = 0.4860776313184267
$data[41,0] = 45071; $data[41,1] = 2023; $data[41,2] = -2.311523918755531; $data[41,3] = 2024; $data[41,4] = -0.1671238976421296
$data[42,0] = 45163; $data[42,1] = 2023; $data[42,2] = -2.156362896191677; $data[42,3] = 2024; $data[42,4] = 0.7251362782769055
$data[43,0] = 45254; $data[43,1] = 2023; $data[43,2] = -2.156362896191677; $data[43,3] = 2024; $data[43,4] = -1.437335768580206
$data[44,0] = 45345; $data[44,1] = 2024; $data[44,2] = -1.379035600217182; $data[44,3] = 2025; $data[44,4] = -1.053916016632561
$data[45,0] = 45436; $data[45,1] = 2024; $data[45,2] = -0.8195740704359578; $data[45,3] = 2025; $data[45,4] = -0.4351858173977874
$data[46,0] = 45534; $data[46,1] = 2024; $data[46,2] = -0.8205034771073372; $data[46,3] = 2025; $data[46,4] = -0.1639718104104904
$data[47,0] = 45618; $data[47,1] = 2024; $data[47,2] = -0.8205034771073372; $data[47,3] = 2025; $data[47,4] = 0.07297157746815053
$data[48,0] = 45713; $data[48,1] = 2025; $data[48,2] = 0.2518138058242414; $data[48,3] = 2026; $data[48,4] = -1.185914719184011
$data[49,0] = 45800; $data[49,1] = 2025; $data[49,2] = 0.4857399523052974; $data[49,3] = 2026; $data[49,4] = -0.3387208997876479
$data[50,0] = 45891; $data[50,1] = 2025; $data[50,2] = 0.5750555200350504; $data[50,3] = 2026; $data[50,4] = 0.3466793681147884

$ws.Range("A2:E52").Value = $data

